$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Type question5..question10 into A6:A11 first
$ws.Cells.Item(6, 1).Value = "question5"
$ws.Cells.Item(7, 1).Value = "question6"
$ws.Cells.Item(8, 1).Value = "question7"
$ws.Cells.Item(9, 1).Value = "question8"
$ws.Cells.Item(10, 1).Value = "question9"
$ws.Cells.Item(11, 1).Value = "question10"

# Step 2: Go back and fill in the missing question4 in A5
$ws.Cells.Item(5, 1).Value = "question4"

# Step 3: Fill the answers column top to bottom
$ws.Cells.Item(5, 2).Value = "answer4"
$ws.Cells.Item(6, 2).Value = "answer5"
$ws.Cells.Item(7, 2).Value = "answer6"
$ws.Cells.Item(8, 2).Value = "answer7"
$ws.Cells.Item(9, 2).Value = "answer8"
$ws.Cells.Item(10, 2).Value = "answer9"
$ws.Cells.Item(11, 2).Value = "answer10"

# Update selection to next empty cell in column B
$ws.Range("B12").Select()
